$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns("D:D").Insert()

Write-Host "Col D width:" $ws.Columns("D").ColumnWidth
Write-Host "Col K width:" $ws.Columns("K").ColumnWidth
Write-Host "Col L width:" $ws.Columns("L").ColumnWidth
